$d = $word.ActiveDocument

# Locate the paragraph that ends with:
#   "git restore --staged [filename] - To remove a file from staging area"
# and insert the new content right after it (and before the trailing blank
# " " paragraph / sectPr that currently follows it).
$target = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $cand = $d.Paragraphs.Item($i)
    if ($cand.Range.Text -like "*git restore --staged*staging area*") {
        $target = $cand
    }
}

if ($target -eq $null) {
    throw "Could not find the 'git restore --staged' paragraph"
}

$full = $target.Range.Duplicate
# Collapse to a point that sits just inside the paragraph's text (before the
# paragraph mark) so InsertXML appends brand new paragraphs right after this
# one instead of splicing into its own run.
$insertionPoint = $d.Range($full.End - 1, $full.End - 1)

$bC = [char]0x2013   # en dash "-"
$lq = [char]0x201C   # left curly quote
$rq = [char]0x201D   # right curly quote

$p1Text = "git log " + $bC + " To see the history of commits or project changes"
$p2Text = "rm -rf [filename] " + $bC + " To delete a file from directory"
$p6Text = "HOW TO RECOVER THE DELETED FILE or GO BACK TO THE PREVIOUS COMMIT:"
$p7Text = "Copy the hash code of commit from " + $lq + "git log" + $rq

$boldRPr = '<w:rPr><w:b/><w:bCs/><w:sz w:val="32"/><w:szCs w:val="32"/></w:rPr>'
$plainRPr = '<w:rPr><w:sz w:val="32"/><w:szCs w:val="32"/></w:rPr>'

$xml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
    '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
    '<pkg:xmlData>' +
    '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
    '<w:body>' +
        '<w:p><w:pPr>' + $boldRPr + '</w:pPr><w:r>' + $boldRPr + '<w:t>' + $p1Text + '</w:t></w:r></w:p>' +
        '<w:p><w:pPr>' + $boldRPr + '</w:pPr><w:r>' + $boldRPr + '<w:t>' + $p2Text + '</w:t></w:r></w:p>' +
        '<w:p><w:pPr>' + $boldRPr + '</w:pPr></w:p>' +
        '<w:p><w:pPr>' + $boldRPr + '</w:pPr></w:p>' +
        '<w:p><w:pPr>' + $boldRPr + '</w:pPr></w:p>' +
        '<w:p><w:pPr>' + $boldRPr + '</w:pPr><w:r>' + $boldRPr + '<w:lastRenderedPageBreak/><w:t>' + $p6Text + '</w:t></w:r></w:p>' +
        '<w:p><w:pPr>' + $plainRPr + '</w:pPr><w:r>' + $plainRPr + '<w:t>' + $p7Text + '</w:t></w:r></w:p>' +
    '</w:body>' +
    '</w:document>' +
    '</pkg:xmlData>' +
    '</pkg:part>' +
    '</pkg:package>'

$insertionPoint.InsertXML($xml)

Write-Output ("Paragraph count now: " + $d.Paragraphs.Count)
